$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix navbar & banner: correct dates and fill in build percentages

# C2 ("Header" row, Finish Build): was free text "12/13" -> a real date 12/11/2023
$ws.Range("C2").NumberFormat = "m/d;@"
$ws.Range("C2").Value = "12/11/2023"
# D2 (% Build for Header row) -> 90%
$ws.Range("D2").Value = 0.9

# C8 ("Banner" row, Finish Build): corrected from 12/13/2020 to 12/11/2020
$ws.Range("C8").Value = "12/11/2020"
# D8 (% Build for Banner row) -> 90%
$ws.Range("D8").Value = 0.9

# B10 ("Restaurant" row, Start Build) -> 12/11/2023
$ws.Range("B10").Value = "12/11/2023"

# Move the active selection to C10
$ws.Range("C10").Select()
